$d = $word.ActiveDocument

function Normalize-Paragraph([int]$index) {
    # Re-serialize a paragraph's runs into a single run, dropping any
    # <w:proofErr/> markers and run-splits left over from Word's
    # spelling/grammar checker, without altering the visible text.
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null
    $original = $r.Text
    # Force a real content change so the run structure is rebuilt, then
    # restore the exact original text (round-trip) in a second pass.
    $r.Text = $original + [char]1
    $r2 = $d.Paragraphs.Item($index).Range
    $r2.MoveEnd(1, -1) | Out-Null
    $r2.Text = $original
}

# ---------------------------------------------------------------------
# Collapse the runs that were only split apart by spell/grammar-checker
# <w:proofErr/> markers back into single runs (no visible text change).
# ---------------------------------------------------------------------
Normalize-Paragraph 4
Normalize-Paragraph 5
Normalize-Paragraph 25
Normalize-Paragraph 31

# ---------------------------------------------------------------------
# Add the new "impact on delay" discussion before the _GoBack bookmark.
# ---------------------------------------------------------------------
$p34 = $d.Paragraphs.Item(34)
$r34 = $p34.Range
$r34.MoveEnd(1, -1) | Out-Null
$part1 = "New simulation has been performed to study the impact on delay. For each figure on message cost, a new simulation has been performed and studied to study the corresponding impact on delay. In summary, the centralized TED will introduce the most delay as expected. "
$part2 = "SPT will usually has smaller delay than distributed TED but when the event probability is low in the network, distributed TED can still outperform SPT."
$startPos = $r34.Start
$r34.InsertBefore($part1 + $part2)

$splitStart = $startPos + $part1.Length
$splitEnd = $splitStart + $part2.Length
$r34b = $d.Range($splitStart, $splitEnd)
$escaped = $part2 -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>' + $escaped + '</w:t></w:r></w:p>'
$r34b.InsertXML($frag)

# ---------------------------------------------------------------------
# Remove the (now stale) lastRenderedPageBreak from the "Section 7.1 has
# been modified..." paragraph - it moved to the new text above.
# ---------------------------------------------------------------------
Normalize-Paragraph 36
